# Delete entire row 3 ("affiliate marketing" / "affiliate.marketing.guide"),
# shifting all subsequent rows up by one, and leave the selection on A3
# (the cell that now occupies the position formerly below the deleted row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Delete()

$ws.Range("A3").Select()
